$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.720.13"
$ws.Range("E2").Value = "  +3.70%  "
$ws.Range("D3").Value = "2.266.83"
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.52%  "
$ws.Range("E7").Value = "  +2.76%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.480"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("D15").Value = "2.615.11"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("D17").Value = "2.247.87"
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("E18").Value = "  +3.54%  "
$ws.Range("D19").Value = "41.628.91"
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.42%  "
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("E25").Value = "  +3.63%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("E30").Value = "  -4.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.54%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.87%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  +3.76%  "
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("E37").Value = "  +1.95%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.69%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.105"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.116"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("D43").Value = "2.051.96"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("E45").Value = "  +4.44%  "
$ws.Range("E46").Value = "  +2.61%  "
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("E48").Value = "  +8.33%  "
$ws.Range("E49").Value = "  +4.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.54%  "
$ws.Range("E51").Value = "  +2.61%  "
